$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "clean-up"/closing-day note row (row 52): the roster slot is no
# longer assigned to 妙柳師姐, the shift note "清屯..." is removed, and the
# shift type becomes 休息日 (rest day).
$ws.Range("B52").Value = ""
$ws.Range("C52").Value = "休息日"
$ws.Range("D52").Value = ""

# Reflect the final selection left on the sheet after the edit.
$ws.Range("B52").Select() | Out-Null
